$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Move the "004756968 / DANIELY" row from its old location (row 85,
#    value 55.51) up near the top of the sheet (new row 8, value
#    5118.32), i.e. delete the old row and insert a fresh one above
#    the "005341184 / BRENO" row.
# ---------------------------------------------------------------------
$ws.Rows.Item(85).Delete()

$ws.Rows.Item(8).Insert()
$ws.Cells.Item(8, 1).NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = "004756968"
$ws.Cells.Item(8, 1).ClearFormats()
$ws.Cells.Item(8, 2).Value = "DANIELY"
$ws.Cells.Item(8, 3).Value = 5118.32

# ---------------------------------------------------------------------
# 2) Add two brand-new rows near the bottom of the data, right after
#    "004400000 / VILMA" (row 192) and before "004387250 / MONICA"
#    (row 193).
# ---------------------------------------------------------------------
$ws.Rows.Item(193).Insert()
$ws.Cells.Item(193, 1).NumberFormat = "@"
$ws.Cells.Item(193, 1).Value = "004644524"
$ws.Cells.Item(193, 1).ClearFormats()
$ws.Cells.Item(193, 2).Value = "PAULO"
$ws.Cells.Item(193, 3).Value = -16.83

$ws.Rows.Item(194).Insert()
$ws.Cells.Item(194, 1).NumberFormat = "@"
$ws.Cells.Item(194, 1).Value = "004222784"
$ws.Cells.Item(194, 1).ClearFormats()
$ws.Cells.Item(194, 2).Value = "RAFAEL"
$ws.Cells.Item(194, 3).Value = -19.41
